$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MapSet")
$ws.Rows.Item(9).Insert()
